$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.933.03"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "1.909.84"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("D4").Value = "'0.9993"
$ws.Range("E4").Value = "  -0.70%  "
$ws.Range("D5").Value = "'313.09"
$ws.Range("E5").Value = "  -0.78%  "
$ws.Range("D6").Value = "'0.9987"
$ws.Range("E6").Value = "  -0.63%  "
$ws.Range("D7").Value = "'0.5010"
$ws.Range("E7").Value = "  +4.14%  "
$ws.Range("D8").Value = "'0.3820"
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").Value = "'0.07323"
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("D10").Value = "'0.9117"
$ws.Range("E10").Value = "  -2.31%  "
$ws.Range("D11").Value = "'21.27"
$ws.Range("E11").Value = "  +2.31%  "
$ws.Range("D12").Value = "'0.07672"
$ws.Range("E12").Value = "  -1.66%  "
$ws.Range("D13").Value = "1.927.60"
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("D14").Value = "'5.488"
$ws.Range("E14").Value = "  -0.21%  "
$ws.Range("D15").Value = "'92.96"
$ws.Range("E15").Value = "  +1.19%  "
$ws.Range("D16").Value = "'0.9994"
$ws.Range("E16").Value = "  -0.75%  "
$ws.Range("D17").Value = "'0.000008748"
$ws.Range("E17").Value = "  -1.48%  "
$ws.Range("D18").Value = "'0.9990"
$ws.Range("E18").Value = "  -0.66%  "
$ws.Range("D19").Value = "27.971.18"
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("D20").Value = "'14.71"
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("D21").Value = "'5.184"
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("D22").Value = "2.170.84"
$ws.Range("E22").Value = "  +1.84%  "
$ws.Range("E23").Value = "  -0.39%  "
$ws.Range("D24").Value = "'6.613"
$ws.Range("E24").Value = "  -0.23%  "
$ws.Range("D25").Value = "'153.19"
$ws.Range("E25").Value = "  -2.15%  "
$ws.Range("D26").Value = "'1.845"
$ws.Range("E26").Value = "  -3.44%  "
$ws.Range("D27").Value = "'2.213"
$ws.Range("E27").Value = "  +3.98%  "
$ws.Range("E28").Value = "  -0.34%  "
$ws.Range("D29").Value = "'115.49"
$ws.Range("E29").Value = "  -1.04%  "
$ws.Range("D30").Value = "'4.932"
$ws.Range("E30").Value = "  -0.64%  "
$ws.Range("D31").Value = "'0.09038"
$ws.Range("E31").Value = "  +1.01%  "
$ws.Range("D32").Value = "'3.208"
$ws.Range("E32").Value = "  -2.73%  "
$ws.Range("D33").Value = "'4.866"
$ws.Range("E33").Value = "  +4.32%  "
$ws.Range("D34").Value = "'1.237"
$ws.Range("E34").Value = "  -1.37%  "
$ws.Range("D35").Value = "'0.7780"
$ws.Range("E35").Value = "  +0.47%  "
$ws.Range("D36").Value = "'0.02088"
$ws.Range("E36").Value = "  +1.95%  "
$ws.Range("D37").Value = "'2.597"
$ws.Range("E37").Value = "  -0.49%  "
$ws.Range("D38").Value = "'3.069"
$ws.Range("E38").Value = "  +2.40%  "
$ws.Range("D39").Value = "'1.094"
$ws.Range("E39").Value = "  -1.32%  "
$ws.Range("D40").Value = "'0.5557"
$ws.Range("E40").Value = "  +0.63%  "
$ws.Range("D41").Value = "'0.05291"
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("D42").Value = "'6.890"
$ws.Range("E42").Value = "  -1.95%  "
$ws.Range("D43").Value = "'113.67"
$ws.Range("E43").Value = "  +4.62%  "
$ws.Range("D44").Value = "'8.531"
$ws.Range("E44").Value = "  +0.37%  "
$ws.Range("D45").Value = "'0.1520"
$ws.Range("E45").Value = "  -0.42%  "
$ws.Range("E46").Value = "  -0.78%  "
$ws.Range("D47").Value = "'0.4842"
$ws.Range("E47").Value = "  +0.28%  "
$ws.Range("D48").Value = "'0.9983"
$ws.Range("E48").Value = "  -0.69%  "
$ws.Range("D49").Value = "'1.643"
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("D50").Value = "'67.71"
$ws.Range("E50").Value = "  -0.31%  "
$ws.Range("D51").Value = "'0.06054"
$ws.Range("E51").Value = "  -0.41%  "
